$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update benchmark timing values (C, D columns) and one new-vars count (B7)
$ws.Range("C2").Value = 0.01676607131958008
$ws.Range("D2").Value = 0.05437588691711426

$ws.Range("C3").Value = 0.06878900527954102
$ws.Range("D3").Value = 0.2041571140289307

$ws.Range("C4").Value = 0.210496187210083
$ws.Range("D4").Value = 0.735598087310791

$ws.Range("C5").Value = 0.4108161926269531
$ws.Range("D5").Value = 1.768811941146851

$ws.Range("C6").Value = 0.7316422462463379
$ws.Range("D6").Value = 4.205510854721069

$ws.Range("B7").Value = 9
$ws.Range("C7").Value = 1.21717095375061
$ws.Range("D7").Value = 11.14751887321472

$ws.Range("C8").Value = 1.702029943466187
$ws.Range("D8").Value = 28.27365374565125

$ws.Range("C9").Value = 2.500485897064209
$ws.Range("D9").Value = 79.25779485702515

$wb.Save()
